$d = $word.ActiveDocument

# 1) Heading text change: "3.1生成DSG实验结果" -> "3.1生成skyline layers实验结果"
$d.Content.Find.Execute("3.1生成DSG实验结果", $true, $false, $false, $false, $false, $true, 1, $false, "3.1生成skyline layers实验结果", 2)

# 2) Trim the long paragraph's wording (two localized removals)
$d.Content.Find.Execute("并进行预处理构造DSG", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
$d.Content.Find.Execute("和构造DSG", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3) Move the "_GoBack" bookmark from the end of the document to right after
#    "通过图表可以发现，" (this is where Word drops it after the last text edit).
$anchor = $d.Content
$anchor.Find.Execute("通过图表可以发现，", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bkRange = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $bkRange)

# 4) Resize the three chart InlineShapes (图表3/4/5) to their new dimensions.
#    InlineShapes 5, 6, 7 correspond to 图表 3, 图表 4, 图表 5 respectively.
$chart3 = $d.InlineShapes.Item(5)
$chart3.Width = 379.5
$chart3.Height = 197.0

$chart4 = $d.InlineShapes.Item(6)
$chart4.Width = 378.05
$chart4.Height = 202.25

$chart5 = $d.InlineShapes.Item(7)
$chart5.Width = 376.85
$chart5.Height = 194.2

Write-Output "done"
